$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.997.62'
$ws.Range("E2").Value = '  -0.39%  '

# Row 3
$ws.Range("D3").Value = '1.805.83'
$ws.Range("E3").Value = '  +0.91%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.26%  '

# Row 5
$ws.Range("E5").Value = '  +0.62%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.28'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.25%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4978'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.78%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3856'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.50%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09327'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +16.66%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.094'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.24%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '40.57'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.99%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.334'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.06%  '

# Row 13
$ws.Range("B13").Value = 'BinanceUSD'
$ws.Range("C13").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.005'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +0.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.70'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.20%  '

# Row 15
$ws.Range("D15").Value = '1.802.33'
$ws.Range("E15").Value = '  +0.84%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.217'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -0.71%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001109'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.69%  '

# Row 18
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.94'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.85%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06573'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.31%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.48%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.11'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.955'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +0.18%  '

# Row 23
$ws.Range("D23").Value = '28.033.99'
$ws.Range("E23").Value = '  -0.40%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.00'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.17%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.232'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.72%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.59'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.19%  '

# Row 27
$ws.Range("D27").Value = '2.018.78'
$ws.Range("E27").Value = '  +1.28%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.45'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.33%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.374'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.16%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.01'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.53%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1075'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.045'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.46%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.544'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +0.23%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.627'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -1.07%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06841'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.63%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '8.860'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +3.57%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02301'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.24%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2133'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.35%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '11.35'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -6.42%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.934'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.51%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6128'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.30%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.145'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.30%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.03'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -1.50%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5870'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.03%  '

# Row 46
$ws.Range("B46").Value = 'WEMIXTOKEN'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.283'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.10%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.663'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -2.63%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.03'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.70%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.944'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +1.48%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.170'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.67%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06731'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.27%  '

